$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "TwoxTwo-ScalarOutNest"

$arr = New-Object 'object[,]' 31,7
$arr[0,1] = "'benchmark"
$arr[0,2] = "'end=1.1"
$arr[0,3] = "'PX=1"
$arr[0,4] = "'PL=1"
$arr[0,5] = "'Itax=0.1"
$arr[0,6] = "'Otax=0.1"
$arr[1,0] = "'X.L"
$arr[1,1] = 1
$arr[1,2] = 1.0087818924459078
$arr[1,3] = 1.0087818984480226
$arr[1,4] = 1.0087818924469572
$arr[1,5] = 0.80862897655434263
$arr[1,6] = 0.97382882042375429
$arr[2,0] = "'Y.L"
$arr[2,1] = 1
$arr[2,2] = 1.0915655837604266
$arr[2,3] = 1.0915655804126048
$arr[2,4] = 1.0915655837428009
$arr[2,5] = 1.2594034491556683
$arr[2,6] = 1.121135937191029
$arr[3,0] = "'U.L"
$arr[3,1] = 1
$arr[3,2] = 1.0503441611049424
$arr[3,3] = 1.0503441602079413
$arr[3,4] = 1.0503441610960742
$arr[3,5] = 1.0533658697164114
$arr[3,6] = 1.0514550815800936
$arr[4,0] = "'XP.L"
$arr[4,1] = 1
$arr[4,2] = 1.0363764562165405
$arr[4,3] = 1.0363764493185874
$arr[4,4] = 1.0363764562119133
$arr[4,5] = 0.958887134059208
$arr[4,6] = 1.022931192686233
$arr[5,0] = "'FR.L"
$arr[5,1] = 1
$arr[5,2] = 1
$arr[5,3] = 1
$arr[5,4] = 1
$arr[5,5] = 1
$arr[5,6] = 1
$arr[6,0] = "'XU.L"
$arr[6,1] = 1
$arr[6,2] = 1.0363764562158884
$arr[6,3] = 1.0363764467002947
$arr[6,4] = 1.03637645621167
$arr[6,5] = 0.95888713402749703
$arr[6,6] = 1.0229311926913409
$arr[7,0] = "'PX.L"
$arr[7,1] = 1
$arr[7,2] = 1.017531272749113
$arr[7,3] = 1
$arr[7,4] = 1.0460796198843898
$arr[7,5] = 1.1303175404901751
$arr[7,6] = 1.1348010195250842
$arr[8,0] = "'PY.L"
$arr[8,1] = 1
$arr[8,2] = 0.98096444410281769
$arr[8,3] = 0.9640631895841022
$arr[8,4] = 1.0084868547104582
$arr[8,5] = 0.92696026052833602
$arr[8,6] = 1.063341999589474
$arr[9,0] = "'PU.L"
$arr[9,1] = 1
$arr[9,2] = 1.0039999207779782
$arr[9,3] = 0.98670178269616649
$arr[9,4] = 1.0321686257962945
$arr[9,5] = 1.0289368379244177
$arr[9,6] = 1.1040161208535313
$arr[10,0] = "'PL.L"
$arr[10,1] = 1
$arr[10,2] = 0.97270920241681125
$arr[10,3] = 0.95595017804950067
$arr[10,4] = 1
$arr[10,5] = 1
$arr[10,6] = 1
$arr[11,0] = "'PK.L"
$arr[11,1] = 1
$arr[11,2] = 1.0360238527648151
$arr[11,3] = 1.0181739667949257
$arr[11,4] = 1.0650910366187067
$arr[11,5] = 0.98436152085413908
$arr[11,6] = 1.050031164244674
$arr[12,0] = "'PF.L"
$arr[12,1] = 1
$arr[12,2] = 1.0545454545441773
$arr[12,3] = 1.0363764440814733
$arr[12,4] = 1.0841322893692942
$arr[12,5] = 1.083846947285189
$arr[12,6] = 1.1608233605436402
$arr[13,0] = "'PXD.L"
$arr[13,1] = 1
$arr[13,2] = 1.0175312727497825
$arr[13,3] = 1.0000000024997451
$arr[13,4] = 1.0460796198846283
$arr[13,5] = 1.1303175405730492
$arr[13,6] = 1.134801019517711
$arr[14,0] = "'PXX.L"
$arr[14,1] = 1
$arr[14,2] = 1.0175312727484433
$arr[14,3] = 0.99999999750029722
$arr[14,4] = 1.0460796198841524
$arr[14,5] = 1.1303175404116717
$arr[14,6] = 1.1348010195324607
$arr[15,0] = "'PKX.L"
$arr[15,1] = 1
$arr[15,2] = 1.0545454545441773
$arr[15,3] = 1.0363764440814733
$arr[15,4] = 1.0841322893692942
$arr[15,5] = 1.083846947285189
$arr[15,6] = 1.1608233605436402
$arr[16,0] = "'SX.L"
$arr[16,1] = 80
$arr[16,2] = 80
$arr[16,3] = 80
$arr[16,4] = 80
$arr[16,5] = 80
$arr[16,6] = 80
$arr[17,0] = "'SXY.L"
$arr[17,1] = 20
$arr[17,5] = 20
$arr[17,6] = 20
$arr[18,0] = "'SY.L"
$arr[18,1] = 80
$arr[18,2] = 80
$arr[18,3] = 80
$arr[18,4] = 80
$arr[18,5] = 80
$arr[18,6] = 80
$arr[19,0] = "'SYX.L"
$arr[19,1] = 40
$arr[19,5] = 40
$arr[19,6] = 40
$arr[20,0] = "'SU.L"
$arr[20,1] = 220
$arr[20,2] = 220
$arr[20,3] = 220
$arr[20,4] = 220
$arr[20,5] = 220
$arr[20,6] = 220
$arr[21,0] = "'DXL.L"
$arr[21,1] = 40
$arr[21,2] = 41.542442674765098
$arr[21,3] = 41.542442691088624
$arr[21,4] = 41.542442673943768
$arr[21,5] = 39.623493954449998
$arr[21,6] = 41.189005357666694
$arr[22,0] = "'DXK.L"
$arr[22,1] = 60
$arr[22,2] = 58.505481567020347
$arr[22,3] = 58.505481551694402
$arr[22,4] = 58.50548156779147
$arr[22,5] = 60.379484236749228
$arr[22,6] = 58.839689849532405
$arr[23,0] = "'DYL.L"
$arr[23,1] = 80
$arr[23,2] = 82.535339520565458
$arr[23,3] = 82.535339547026339
$arr[23,4] = 82.5353395192341
$arr[23,5] = 79.370351203301908
$arr[23,6] = 81.960586873948159
$arr[24,0] = "'DYK.L"
$arr[24,1] = 40
$arr[24,2] = 37.543075921873587
$arr[24,3] = 37.54307589702983
$arr[24,4] = 37.54307592312361
$arr[24,5] = 40.634635036712858
$arr[24,6] = 38.08656279327819
$arr[25,0] = "'DUX.L"
$arr[25,1] = 80
$arr[25,2] = 78.936142616180291
$arr[25,3] = 78.936142557576616
$arr[25,4] = 78.936142616709958
$arr[25,5] = 72.82462143888138
$arr[25,6] = 77.829758828574981
$arr[26,0] = "'DUY.L"
$arr[26,1] = 100
$arr[26,2] = 102.34824787101159
$arr[26,3] = 102.34824777361034
$arr[26,4] = 102.34824786999235
$arr[26,5] = 111.00118108759018
$arr[26,6] = 103.82512129465007
$arr[27,0] = "'DU.L"
$arr[27,1] = 220
$arr[27,2] = 231.07571544002587
$arr[27,3] = 231.07571525715505
$arr[27,4] = 231.07571544123067
$arr[27,5] = 231.74049138660823
$arr[27,6] = 231.32011801351359
$arr[28,0] = "'CWI.L"
$arr[28,1] = 220
$arr[28,2] = 231.07571544002587
$arr[28,3] = 231.07571525715505
$arr[28,4] = 231.07571544123067
$arr[28,5] = 231.74049138660823
$arr[28,6] = 231.32011801351359
$arr[29,0] = "'RA.L"
$arr[29,1] = 220
$arr[29,2] = 232
$arr[29,3] = 228.00282018202665
$arr[29,4] = 238.50910366187065
$arr[29,5] = 238.44632842638742
$arr[29,6] = 255.38113936466033
$arr[30,0] = "'ROW.L"
$arr[30,1] = 20
$arr[30,2] = 21.090909090883546
$arr[30,3] = 20.727528881629468
$arr[30,4] = 21.682645787385884
$arr[30,5] = 21.676938945703782
$arr[30,6] = 23.216467210872803

$ws.Range("A1:G31").Value = $arr
